# Commit: "Allowing resource optimized OpenStack"
#
# One option for the deployment of edge cloud infrastructures is to use
# different instances of the same resource optimized OpenStack in both the
# large, medium and small edge sites. This drops the "LW" (lightweight)
# qualifier from the label of the OpenStack box in the figure, since the
# box now simply reads "OpenStack".
#
# The text lives on slide 1 (sldId 290), inside shape id 16 ("圆角矩形 13",
# a rounded-rectangle label), which itself is nested one level deep inside
# the single top-level group shape ("组合 30") that makes up the whole
# figure on that slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The entire figure on this slide is one big group shape.
$figureGroup = $s.Shapes.Item(1)
$groupItems = $figureGroup.GroupItems

$target = $null
for ($i = 1; $i -le $groupItems.Count; $i++) {
    $shp = $groupItems.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "LW OpenStack") {
                $target = $shp
            }
        }
    }
}

if ($target -ne $null) {
    $target.TextFrame.TextRange.Text = "OpenStack"
} else {
    # Fallback: the label shape is known to be shape id 16 / index 8 in the
    # group's GroupItems collection.
    $fallback = $groupItems.Item(8)
    if ($fallback.TextFrame.TextRange.Text -eq "LW OpenStack") {
        $fallback.TextFrame.TextRange.Text = "OpenStack"
    }
}
